# Updated cryptos list on Sat Apr 20 07:26:34 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the coin rows with
# the latest scraped values. One pair of rows (Uniswap / InternetComputer)
# also changed rank position, so their Coin/Link/Price/Volume cells are
# rewritten in place to reflect the new ordering.
#
# D-column values are plain text (e.g. "64.020.77", using dots as
# thousands separators) in the source data, so NumberFormat is forced to
# "@" (Text) before assigning values that would otherwise be
# auto-recognized as numbers (e.g. "143.76"), keeping them as literal
# text strings like the original sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.020.77"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.063.55"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.76"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.062.75"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("E11").Value = "  -10.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.498"
$ws.Range("E12").Value = "  +9.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.84"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.565.14"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.080.01"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.065.53"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.81"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.61"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.35"
$ws.Range("E24").Value = "  +8.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.71"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.35"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.47"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0413"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "452.22"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0816"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.033.52"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.75"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +7.22%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.14"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  +1.14%  "